$wb = $excel.ActiveWorkbook

# --- Sheet "registrants" ---
$ws1 = $wb.Worksheets.Item("registrants")

# D7 displayed value changes from "PT_BR" to "IT_IT"
$ws1.Range("D7").Value = "IT_IT"

# --- Sheet "services" ---
$ws2 = $wb.Worksheets.Item("services")

# Update Language column values
$ws2.Range("C3").Value = "DE_DE"
$ws2.Range("C5").Value = "EN_GB"
$ws2.Range("C6").Value = "EN_GB"
$ws2.Range("C7").Value = "FR_FR"

# (B8 / C8 keep their original displayed text "NEWS" / "DE_DE" - no edit
# needed there, only their underlying shared-string slot shuffles around,
# which is an internal/non-semantic storage detail.)

# Selection on "services" changes to range A2:D8 with active cell A2
$ws2.Range("A2:D8").Select()

# Selection on "registrants" changes to active cell A3. Re-selecting here
# (last) also re-activates "registrants" as the selected tab, matching the
# original workbook state (tabSelected stays on "registrants").
$ws1.Range("A3").Select()
